# Swap the taxon-observation data between row 2 and row 3.
# (Columns A, B, D, E, F, G, H, Q, R hold values that need to trade places;
#  the other columns already contain identical data in both rows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $r2 = $col + "2"
    $r3 = $col + "3"
    $v2 = $ws.Range($r2).Value()
    $v3 = $ws.Range($r3).Value()
    $ws.Range($r2).Value = $v3
    $ws.Range($r3).Value = $v2
}

# The empty placeholder cell in column L moves from row 3 to row 2.
$ws.Range("L3").ClearContents()
$ws.Range("L2").Value = ""
